$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.921.13'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.875.10'
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''0.7390'
$ws.Range("E5").Value = '  -4.65%  '
$ws.Range("D6").Value = '''242.59'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''0.3154'
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").Value = '''0.07154'
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").Value = '''24.68'
$ws.Range("E10").Value = '  -4.23%  '
$ws.Range("D11").Value = '''0.08408'
$ws.Range("E11").Value = '  -3.27%  '
$ws.Range("D12").Value = '''0.7507'
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").Value = '''5.426'
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = '1.885.93'
$ws.Range("E14").Value = '  -5.63%  '
$ws.Range("D15").Value = '''92.60'
$ws.Range("E15").Value = '  -1.98%  '
$ws.Range("D16").Value = '29.896.00'
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '''6.097'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").Value = '''13.59'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").Value = '''242.99'
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D20").Value = '''0.000007818'
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Value = '''0.9995'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '2.115.25'
$ws.Range("E22").Value = '  -6.64%  '
$ws.Range("D23").Value = '''7.991'
$ws.Range("E23").Value = '  -2.50%  '
$ws.Range("D24").Value = '''0.9998'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '''0.1550'
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").Value = '''9.296'
$ws.Range("D27").Value = '''165.02'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").Value = '''18.61'
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("D29").Value = '''2.036'
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  +4.15%  '
$ws.Range("D31").Value = '''4.595'
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("D32").Value = '''1.528'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").Value = '''4.251'
$ws.Range("E33").Value = '  +2.95%  '
$ws.Range("D34").Value = '''0.05321'
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").Value = '''0.7545'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '''0.9957'
$ws.Range("E37").Value = '  -0.61%  '
$ws.Range("D38").Value = '''2.696'
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("D41").Value = '''0.4506'
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("D42").Value = '1.110.61'
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("D43").Value = '''6.048'
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").Value = '''72.21'
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").Value = '''0.8561'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").Value = '''1.002'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").Value = '''103.17'
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").Value = '''7.650'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = '''3.088'
$ws.Range("E49").Value = '  +3.09%  '
$ws.Range("D50").Value = '''1.840'
$ws.Range("E50").Value = '  -2.48%  '
$ws.Range("D51").Value = '2.014.34'
$ws.Range("E51").Value = '  -9.63%  '
